$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: A1 was a shared string "Example"; now becomes numeric 100
$ws1.Range("A1").Value = 100

# Sheet1: add rows 28-30 with values 27, 28, 29
$ws1.Range("A28").Value = 27
$ws1.Range("A29").Value = 28
$ws1.Range("A30").Value = 29

# Update selections
$ws1.Range("K10").Select()
$ws2.Range("A18").Select()

# Activate Sheet1 so it becomes the active tab
$ws1.Activate()
